# Weekly data refresh: insert a new observation as the newest row in the
# "Terminal La Palmera de La Serena - Cebollín" table. The table is sorted
# with the newest date first, so the new row lands at row 200 (just below
# the header and the still-newer rows above it), pushing the previous
# rows 200-222 down to 201-223.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 200, shifting rows 200:222 down
# to 201:223 (and bumping the sheet dimension from R222 to R223).
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with this week's observation.
$ws.Range("A200").Value = 8
$ws.Range("B200").Value = "Terminal La Palmera de La Serena"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44776
$ws.Range("E200").Value = 4
$ws.Range("F200").Value = 100112037
$ws.Range("G200").Value = "Cebollín"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 1300
$ws.Range("K200").Value = 1400
$ws.Range("L200").Value = 1600
$ws.Range("M200").Value = 1500
$ws.Range("N200").Value = "$/paquete 6 unidades"
$ws.Range("O200").Value = "Provincia del Elquí"
$ws.Range("P200").Value = 250
$ws.Range("Q200").Value = 6
$ws.Range("R200").Value = "Hortaliza"
